$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 1.422753333333333
$ws.Range("H2").Value = 4.26826
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 3.76917
$ws.Range("N2").Value = 11.30751
$ws.Range("O2").Value = 0.02686645020528053
$ws.Range("P2").Value = 0.02686645020528053
$ws.Range("Q2").Value = 5.3625991814
$ws.Range("R2").Value = 48.2633926326
$ws.Range("S2").Value = 0.02686645020528053
$ws.Range("T2").Value = 0.02686645020528053

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 1.422753333333333
$ws.Range("H3").Value = 4.26826
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 82.48060333333333
$ws.Range("N3").Value = 247.44181
$ws.Range("O3").Value = 0.5879175050094569
$ws.Range("P3").Value = 0.5879175050094571
$ws.Range("Q3").Value = 117.3495533278444
$ws.Range("R3").Value = 1056.1459799506
$ws.Range("S3").Value = 0.5879175050094569
$ws.Range("T3").Value = 0.5879175050094571

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 1.422753333333333
$ws.Range("H4").Value = 4.26826
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.5716056666666667
$ws.Range("N4").Value = 1.714817
$ws.Range("O4").Value = 0.004074375838860061
$ws.Range("P4").Value = 0.004074375838860062
$ws.Range("Q4").Value = 0.8132538676022222
$ws.Range("R4").Value = 7.31928480842
$ws.Range("S4").Value = 0.004074375838860061
$ws.Range("T4").Value = 0.004074375838860062

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 1.422753333333333
$ws.Range("H5").Value = 4.26826
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 53.47143866666666
$ws.Range("N5").Value = 160.414316
$ws.Range("O5").Value = 0.3811416689464024
$ws.Range("P5").Value = 0.3811416689464024
$ws.Range("Q5").Value = 76.07666760112889
$ws.Range("R5").Value = 684.6900084101599
$ws.Range("S5").Value = 0.3811416689464024
$ws.Range("T5").Value = 0.3811416689464024
